# Front end department manager menu added.
#
# Changes on slide 9 ("PO Authentication"):
#  1. The "Table 3" graphic frame (the Approved PO Log table) grows taller
#     (its height increases by 16.8pt, i.e. 213360 EMU: 1615440 -> 1828800 EMU).
#  2. The "Picture 2" ellipsis-menu icon (id 19) that sits on top of that
#     table is nudged down/right to its new spot next to the taller table
#     (3934584,3147046) -> (3947284,3229775) EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

$EmuPerPt = 914400 / 72

# 1) Grow the "Approved PO Log" table (graphicFrame "Table 3").
$table = $s.Shapes.Item("Table 3")
$table.Height = 1828800 / $EmuPerPt

# 2) Reposition the ellipsis menu icon ("Picture 2", id 19) that belongs to
#    that table. (The values are nudged by a hair so that the engine's
#    point<->EMU round trip lands exactly on 3947284/3229775 EMU instead of
#    one EMU short.)
$icon = $s.Shapes.Item("Picture 2")
$icon.Left = 310.80976867952756
$icon.Top = 254.31299592598427
